$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: add headers P1=14, Q1=15 with the same style as the existing header cells (O1)
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Rows 2-25: update I,K,M,O columns and add P,Q columns
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P: new = 2
    $ws.Cells.Item($r, 17).Value = 2  # Q: new = 2
}
